# This script applies the milestone-2 scrum report edits: several table
# cells referring to "test cases" are reworded to refer to the "test plan"
# instead (part of renaming the testing deliverable from "Test Cases" to
# "Test Plan" throughout the document).

$d = $word.ActiveDocument

function Replace-InCell($table, $row, $col, $findText, $replaceText) {
    $cellRange = $table.Cell($row, $col).Range
    $cellRange.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1) | Out-Null
}

# Table 3 ("Tasks Completed" scrum table)
Replace-InCell $d.Tables(3) 3 2 "Managing Jira and test cases" "Managing Jira and test plan"
Replace-InCell $d.Tables(3) 4 2 "Committing the changes to GitHub repo and test cases" "Committing the changes to GitHub repo and test plan"
Replace-InCell $d.Tables(3) 6 2 "Making the Test cases" "Updating the Test plan"

# Table 5 (Discussion Summary / Outcome table)
Replace-InCell $d.Tables(5) 4 1 "Test Cases" "Test Plan"
Replace-InCell $d.Tables(5) 4 2 "Discussed about the test cases" "Discussed about the test plan"
Replace-InCell $d.Tables(5) 4 3 "Test cases were developed." "Test plan was updated."

# Table 6 (Decision / Rationale table)
Replace-InCell $d.Tables(6) 2 1 "Test Cases to be done" "Test Plan to be updated"

# Table 7 (Task Attempted table)
Replace-InCell $d.Tables(7) 5 2 "Test Cases" "Test Plan"
Replace-InCell $d.Tables(7) 6 2 "Test Cases and Jira management" "Test Plan and Jira management"
Replace-InCell $d.Tables(7) 7 2 "Test cases and git management" "Test Plan and GitHub management"

# Table 8 (Group Member / Task Description table)
Replace-InCell $d.Tables(8) 6 2 "Working on the test cases" "Working on the test plan"
